# Generate Report for Handoff
# Replaces the localization-status report rows for the file that was handed
# off (54668c75-76de-4c0a-95ca-b600d9003c0e.md -> 1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md)
# and the file that is now ready for handoff
# (9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md -> ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md)
# across the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview"
# -----------------------------------------------------------------
$ov = $wb.Worksheets.Item(1)

$ov.Range("A2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$ov.Range("B2").Value = "e2e\1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-15 20:58:41"

$ov.Range("A3").Value = "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$ov.Range("B3").Value = "e2e\ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-15 20:58:41"

# Hyperlinks.Delete() on any range wipes the whole-sheet hyperlink
# collection in this runtime, so wipe then rebuild the two we need,
# keeping the original link targets intact.
$ov.Range("A1").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/54668c75-76de-4c0a-95ca-b600d9003c0e.md", "", "", "e2e\1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md", "", "", "e2e\ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md")

$ov.Columns.Item(5).ColumnWidth = 16.333333333333332
$ov.Columns.Item(6).ColumnWidth = 16.333333333333332

# -----------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------
$zh = $wb.Worksheets.Item(2)

$zh.Range("A2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("G2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-15 20:58:36"
$zh.Range("I2").Style = "Normal"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "'True"
$zh.Range("G3").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-15 20:58:36"
$zh.Range("I3").Style = "Normal"
$zh.Range("I3").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"

$zh.Range("A1").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/54668c75-76de-4c0a-95ca-b600d9003c0e.md", "", "", "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md", "", "", "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md")

$zh.Columns.Item(3).ColumnWidth = 16.333333333333332
$zh.Columns.Item(9).ColumnWidth = 17.833333333333332
$zh.Columns.Item(10).ColumnWidth = 20.833333333333332

# -----------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------
$de = $wb.Worksheets.Item(3)

$de.Range("A2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("G2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.de-de.xlf"
$de.Range("H2").Value = "2016-08-15 20:58:41"
$de.Range("I2").Style = "Normal"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Range("A3").Value = "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "'True"
$de.Range("G3").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.de-de.xlf"
$de.Range("H3").Value = "2016-08-15 20:58:41"
$de.Range("I3").Style = "Normal"
$de.Range("I3").Value = ""
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"

$de.Range("A1").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/54668c75-76de-4c0a-95ca-b600d9003c0e.md", "", "", "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md", "", "", "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md")

$de.Columns.Item(3).ColumnWidth = 16.333333333333332
$de.Columns.Item(9).ColumnWidth = 17.833333333333332
$de.Columns.Item(10).ColumnWidth = 20.833333333333332

$wb.Save()
